$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 3).Value = 19.0
$ws.Cells.Item(2, 4).Value = 8.47
$ws.Cells.Item(2, 5).Value = 1.26
$ws.Cells.Item(2, 6).Value = 5.53
$ws.Cells.Item(2, 7).Value = 1.07
$ws.Cells.Item(3, 3).Value = 19.0
$ws.Cells.Item(3, 4).Value = 5.0
$ws.Cells.Item(3, 5).Value = 1.76
$ws.Cells.Item(3, 6).Value = 4.21
$ws.Cells.Item(3, 7).Value = 0.63
$ws.Cells.Item(4, 3).Value = 38.0
$ws.Cells.Item(4, 4).Value = 7.42
$ws.Cells.Item(4, 5).Value = 1.84
$ws.Cells.Item(4, 6).Value = 4.89
$ws.Cells.Item(5, 3).Value = 19.0
$ws.Cells.Item(5, 4).Value = 7.79
$ws.Cells.Item(5, 5).Value = 1.78
$ws.Cells.Item(5, 6).Value = 5.42
$ws.Cells.Item(5, 7).Value = 1.22
$ws.Cells.Item(6, 3).Value = 19.0
$ws.Cells.Item(6, 4).Value = 7.74
$ws.Cells.Item(6, 5).Value = 1.79
$ws.Cells.Item(6, 6).Value = 5.21
$ws.Cells.Item(6, 7).Value = 1.27
$ws.Cells.Item(7, 3).Value = 19.0
$ws.Cells.Item(7, 4).Value = 6.58
$ws.Cells.Item(7, 5).Value = 1.87
$ws.Cells.Item(7, 6).Value = 4.32
$ws.Cells.Item(7, 7).Value = 1.06
$ws.Cells.Item(8, 3).Value = 19.0
$ws.Cells.Item(8, 4).Value = 6.63
$ws.Cells.Item(8, 5).Value = 1.98
$ws.Cells.Item(8, 6).Value = 4.47
$ws.Cells.Item(8, 7).Value = 0.84
$ws.Cells.Item(9, 3).Value = 19.0
$ws.Cells.Item(9, 4).Value = 8.21
$ws.Cells.Item(9, 5).Value = 1.47
$ws.Cells.Item(9, 6).Value = 5.58
$ws.Cells.Item(9, 7).Value = 0.96
$ws.Cells.Item(10, 3).Value = 19.0
$ws.Cells.Item(10, 4).Value = 7.42
$ws.Cells.Item(10, 5).Value = 1.92
$ws.Cells.Item(10, 6).Value = 5.47
$ws.Cells.Item(10, 7).Value = 1.35
$ws.Cells.Item(11, 3).Value = 19.0
$ws.Cells.Item(11, 4).Value = 6.58
$ws.Cells.Item(11, 5).Value = 2.29
$ws.Cells.Item(11, 6).Value = 4.53
$ws.Cells.Item(11, 7).Value = 1.02
$ws.Cells.Item(12, 3).Value = 19.0
$ws.Cells.Item(12, 4).Value = 8.05
$ws.Cells.Item(12, 5).Value = 0.85
$ws.Cells.Item(12, 6).Value = 5.47
$ws.Cells.Item(12, 7).Value = 1.12
$ws.Cells.Item(13, 3).Value = 19.0
$ws.Cells.Item(13, 4).Value = 8.11
$ws.Cells.Item(13, 5).Value = 1.66
$ws.Cells.Item(13, 6).Value = 6.05
$ws.Cells.Item(13, 7).Value = 0.97
$ws.Cells.Item(14, 3).Value = 19.0
$ws.Cells.Item(14, 4).Value = 8.42
$ws.Cells.Item(14, 5).Value = 1.22
$ws.Cells.Item(14, 7).Value = 1.56
$ws.Cells.Item(15, 3).Value = 19.0
$ws.Cells.Item(15, 4).Value = 6.95
$ws.Cells.Item(15, 5).Value = 2.07
$ws.Cells.Item(15, 6).Value = 5.05
$ws.Cells.Item(15, 7).Value = 1.08
$ws.Cells.Item(16, 3).Value = 19.0
$ws.Cells.Item(16, 4).Value = 7.58
$ws.Cells.Item(16, 5).Value = 2.14
$ws.Cells.Item(16, 6).Value = 4.74
$ws.Cells.Item(16, 7).Value = 1.66
$ws.Cells.Item(17, 3).Value = 19.0
$ws.Cells.Item(17, 4).Value = 6.84
$ws.Cells.Item(17, 5).Value = 1.74
$ws.Cells.Item(17, 6).Value = 5.16
$ws.Cells.Item(17, 7).Value = 1.12
$ws.Cells.Item(18, 3).Value = 19.0
$ws.Cells.Item(18, 4).Value = 5.47
$ws.Cells.Item(18, 5).Value = 2.29
$ws.Cells.Item(18, 6).Value = 4.47
$ws.Cells.Item(18, 7).Value = 0.9
$ws.Cells.Item(19, 3).Value = 19.0
$ws.Cells.Item(19, 4).Value = 7.26
$ws.Cells.Item(19, 5).Value = 2.1
$ws.Cells.Item(19, 6).Value = 5.26
$ws.Cells.Item(19, 7).Value = 1.37
$ws.Cells.Item(20, 3).Value = 19.0
$ws.Cells.Item(20, 4).Value = 5.53
$ws.Cells.Item(20, 5).Value = 2.34
$ws.Cells.Item(20, 6).Value = 4.37
$ws.Cells.Item(20, 7).Value = 0.83
$ws.Cells.Item(21, 3).Value = 19.0
$ws.Cells.Item(21, 4).Value = 6.37
$ws.Cells.Item(21, 5).Value = 2.45
$ws.Cells.Item(21, 6).Value = 4.37
$ws.Cells.Item(21, 7).Value = 1.21
$ws.Cells.Item(22, 3).Value = 19.0
$ws.Cells.Item(22, 4).Value = 7.68
$ws.Cells.Item(22, 5).Value = 2.11
$ws.Cells.Item(22, 6).Value = 5.42
$ws.Cells.Item(22, 7).Value = 1.22
$ws.Cells.Item(23, 3).Value = 19.0
$ws.Cells.Item(23, 4).Value = 7.89
$ws.Cells.Item(23, 5).Value = 1.33
$ws.Cells.Item(23, 6).Value = 5.16
$ws.Cells.Item(23, 7).Value = 0.83
$ws.Cells.Item(24, 3).Value = 19.0
$ws.Cells.Item(24, 4).Value = 8.26
$ws.Cells.Item(24, 6).Value = 5.16
$ws.Cells.Item(24, 7).Value = 1.07
$ws.Cells.Item(25, 3).Value = 19.0
$ws.Cells.Item(25, 4).Value = 7.37
$ws.Cells.Item(25, 5).Value = 1.57
$ws.Cells.Item(25, 6).Value = 5.16
$ws.Cells.Item(25, 7).Value = 1.12
$ws.Cells.Item(26, 3).Value = 19.0
$ws.Cells.Item(26, 4).Value = 7.47
$ws.Cells.Item(26, 5).Value = 1.87
$ws.Cells.Item(26, 6).Value = 5.58
$ws.Cells.Item(26, 7).Value = 1.17
$ws.Cells.Item(27, 3).Value = 19.0
$ws.Cells.Item(27, 4).Value = 7.68
$ws.Cells.Item(27, 5).Value = 1.97
$ws.Cells.Item(27, 6).Value = 5.47
$ws.Cells.Item(27, 7).Value = 1.26
$ws.Cells.Item(28, 3).Value = 19.0
$ws.Cells.Item(28, 4).Value = 5.26
$ws.Cells.Item(28, 5).Value = 2.21
$ws.Cells.Item(28, 6).Value = 4.32
$ws.Cells.Item(28, 7).Value = 0.75
$ws.Cells.Item(29, 3).Value = 19.0
$ws.Cells.Item(29, 4).Value = 8.11
$ws.Cells.Item(29, 5).Value = 1.33
$ws.Cells.Item(29, 6).Value = 5.32
$ws.Cells.Item(29, 7).Value = 1.06
$ws.Cells.Item(30, 3).Value = 19.0
$ws.Cells.Item(30, 4).Value = 8.53
$ws.Cells.Item(30, 5).Value = 1.71
$ws.Cells.Item(30, 6).Value = 6.0
$ws.Cells.Item(30, 7).Value = 1.33
$ws.Cells.Item(31, 3).Value = 19.0
$ws.Cells.Item(31, 4).Value = 6.68
$ws.Cells.Item(31, 5).Value = 2.58
$ws.Cells.Item(31, 6).Value = 4.32
$ws.Cells.Item(31, 7).Value = 0.82
$ws.Cells.Item(32, 3).Value = 19.0
$ws.Cells.Item(32, 4).Value = 6.05
$ws.Cells.Item(32, 5).Value = 1.96
$ws.Cells.Item(32, 6).Value = 4.79
$ws.Cells.Item(32, 7).Value = 0.98
